# Observation_periods.xlsx edit
# - Rename header "Interface appearing" -> "Last image with full ROI"
# - Tweak several ROI end-column (D) values (moving ROI for the fall)
# - Populate row 7 (F_h2_f1000_1_s) C/D cells with "nan" placeholders
# - Move active selection to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename column D header text.
$ws.Range("D1").Value = "Last image with full ROI"

# Updated ROI values in column D (end-of-interface / ROI boundary values).
$ws.Range("D2").Value = 1564
$ws.Range("D3").Value = 696
$ws.Range("D4").Value = 709
$ws.Range("D5").Value = 893
$ws.Range("D6").Value = 308
$ws.Range("D8").Value = 377
$ws.Range("D10").Value = 582
$ws.Range("D11").Value = 819

# Row 7 (F_h2_f1000_1_s) previously had no C/D values -- now marked "nan".
$ws.Range("C7").Value = "nan"
$ws.Range("D7").Value = "nan"

# Restore the active cell selection used when the file was last saved.
$ws.Range("E6").Select()
